$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

$ws.Range("H28").Value = 334.33334
$ws.Range("I28").Value = 339.2143
$ws.Range("K28").Value = 339.2143
$ws.Range("M28").Value = 145.7857

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H111").Value = 250
$ws.Range("I111").Value = 250
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 750
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 2317
$ws.Range("N111").ClearContents()

$ws.Range("H112").Value = 4948.5
$ws.Range("J112").Value = 4948.5
$ws.Range("L112").Value = 14845.5
$ws.Range("N112").Value = -17061.5

$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442

$ws.Range("H137").Value = 1887.8125
$ws.Range("I137").Value = 1796.4166
$ws.Range("K137").Value = 5389.2498
$ws.Range("M137").Value = -2839.2498

$ws.Range("H138").Value = 3126.6897
$ws.Range("I138").Value = 1221.4445
$ws.Range("K138").Value = 3664.3335
$ws.Range("M138").Value = 1475.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4241.423
$ws.Range("I32").Value = 3282.261
$ws.Range("K32").Value = 3282.261
$ws.Range("M32").Value = -2995.261

$ws.Range("H61").Value = 3733.111
$ws.Range("I61").Value = 3324.875
$ws.Range("J61").Value = 6999
$ws.Range("K61").Value = 3324.875
$ws.Range("L61").Value = 6999
$ws.Range("M61").Value = -3112.875
$ws.Range("N61").Value = -7423

$ws.Range("H74").Value = 1862.0769
$ws.Range("I74").Value = 1603.1666
$ws.Range("K74").Value = 1603.1666
$ws.Range("M74").Value = -729.1666

$ws.Range("H77").Value = 1862.0769
$ws.Range("I77").Value = 1603.1666
$ws.Range("K77").Value = 8015.833000000001
$ws.Range("M77").Value = -3647.833000000001

$ws.Range("H88").Value = 605.1429000000001
$ws.Range("J88").Value = 797.3333
$ws.Range("L88").Value = 797.3333
$ws.Range("N88").Value = -1609.3333

$ws.Range("H91").Value = 605.1429000000001
$ws.Range("J91").Value = 797.3333
$ws.Range("L91").Value = 797.3333
$ws.Range("N91").Value = -3605.3333

$ws.Range("H132").Value = 3589.8
$ws.Range("I132").Value = 3245
$ws.Range("K132").Value = 9735
$ws.Range("M132").Value = -7205

$ws.Range("H136").Value = 3733.111
$ws.Range("I136").Value = 3324.875
$ws.Range("J136").Value = 6999
$ws.Range("K136").Value = 9974.625
$ws.Range("L136").Value = 20997
$ws.Range("M136").Value = -7424.625
$ws.Range("N136").Value = -26097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 11630.4
$ws.Range("I82").Value = 11630.4
$ws.Range("K82").Value = 11630.4
$ws.Range("M82").Value = -11247.4

$ws.Range("H85").Value = 11630.4
$ws.Range("I85").Value = 11630.4
$ws.Range("K85").Value = 11630.4
$ws.Range("M85").Value = -10304.4

$ws.Range("H107").Value = 13172.182
$ws.Range("I107").Value = 5321.6665
$ws.Range("K107").Value = 5321.6665
$ws.Range("M107").Value = -3401.6665

$ws.Range("H134").Value = 10289.143
$ws.Range("I134").Value = 5999.5
$ws.Range("K134").Value = 17998.5
$ws.Range("M134").Value = -15463.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3921.1428
$ws.Range("I58").Value = 4241.3335
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 4241.3335
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -4038.3335
$ws.Range("N58").Value = -2406

$ws.Range("H132").Value = 1988.8889
$ws.Range("I132").Value = 1060
$ws.Range("K132").Value = 3180
$ws.Range("M132").Value = -650

$ws.Range("H134").Value = 2999.75
$ws.Range("I134").Value = 2499.5
$ws.Range("K134").Value = 7498.5
$ws.Range("M134").Value = -4963.5

$ws.Range("H136").Value = 3921.1428
$ws.Range("I136").Value = 4241.3335
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 12724.0005
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -10174.0005
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 728.5714
$ws.Range("I75").Value = 550
$ws.Range("J75").Value = 800
$ws.Range("K75").Value = 1650
$ws.Range("L75").Value = 2400
$ws.Range("M75").Value = -652
$ws.Range("N75").Value = -4396

$ws.Range("H78").Value = 728.5714
$ws.Range("I78").Value = 550
$ws.Range("J78").Value = 800
$ws.Range("K78").Value = 4950
$ws.Range("L78").Value = 7200
$ws.Range("M78").Value = 42
$ws.Range("N78").Value = -17184

$ws.Range("H113").Value = 944.8333
$ws.Range("J113").Value = 930.8889
$ws.Range("L113").Value = 2792.6667
$ws.Range("N113").Value = -7132.6667

$ws.Range("H117").Value = 1373.75
$ws.Range("J117").Value = 1565
$ws.Range("L117").Value = 4695
$ws.Range("N117").Value = -11579

$ws.Range("H131").Value = 1231.625
$ws.Range("J131").Value = 881
$ws.Range("L131").Value = 2643
$ws.Range("N131").Value = -12723

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1231.3334
$ws.Range("I22").Value = 694
$ws.Range("K22").Value = 694
$ws.Range("M22").Value = -399

$ws.Range("H27").Value = 1231.3334
$ws.Range("I27").Value = 694
$ws.Range("K27").Value = 694
$ws.Range("M27").Value = -587

$ws.Range("H46").Value = 823
$ws.Range("I46").Value = 764
$ws.Range("K46").Value = 764
$ws.Range("M46").Value = -576

$ws.Range("H68").Value = 2574.4
$ws.Range("I68").Value = 2574.4
$ws.Range("K68").Value = 2574.4
$ws.Range("M68").Value = -1825.4

$ws.Range("H71").Value = 2574.4
$ws.Range("I71").Value = 2574.4
$ws.Range("K71").Value = 12872
$ws.Range("M71").Value = -9128

$ws.Range("H82").Value = 1277.3334
$ws.Range("I82").Value = 1420.1428
$ws.Range("J82").Value = 777.5
$ws.Range("K82").Value = 1420.1428
$ws.Range("L82").Value = 777.5
$ws.Range("M82").Value = -1059.1428
$ws.Range("N82").Value = -1499.5

$ws.Range("H85").Value = 1277.3334
$ws.Range("I85").Value = 1420.1428
$ws.Range("J85").Value = 777.5
$ws.Range("K85").Value = 1420.1428
$ws.Range("L85").Value = 777.5
$ws.Range("M85").Value = -172.1428000000001
$ws.Range("N85").Value = -3273.5

$ws.Range("H132").Value = 2803.125
$ws.Range("J132").Value = 1987.5
$ws.Range("L132").Value = 5962.5
$ws.Range("N132").Value = -11022.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 37214
$ws.Range("J46").Value = 37214
$ws.Range("L46").Value = 37214
$ws.Range("N46").Value = -37676

$ws.Range("H81").Value = 1432.5555
$ws.Range("I81").Value = 986.625
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 1973.25
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -912.25
$ws.Range("N81").Value = -12122

$ws.Range("H84").Value = 1432.5555
$ws.Range("I84").Value = 986.625
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 9866.25
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -4562.25
$ws.Range("N84").Value = -60608

$ws.Range("H132").Value = 3315
$ws.Range("I132").Value = 2186.238
$ws.Range("J132").Value = 7265.6665
$ws.Range("K132").Value = 6558.714
$ws.Range("L132").Value = 21796.9995
$ws.Range("M132").Value = -4028.714
$ws.Range("N132").Value = -26856.9995

$ws.Range("H134").Value = 37214
$ws.Range("J134").Value = 37214
$ws.Range("L134").Value = 111642
$ws.Range("N134").Value = -116712

$ws.Range("H136").Value = 1835.238
$ws.Range("I136").Value = 1627.05
$ws.Range("K136").Value = 4881.15
$ws.Range("M136").Value = -2331.15
